$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.316.08'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -3.74%  '
# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.930.06'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -4.00%  '
# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.18%  '
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '248.44'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.03%  '
# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.7143'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -7.01%  '
# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.23%  '
# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3255'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -9.73%  '
# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '27.23'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.75%  '
# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06806'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.29%  '
# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8004'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -6.19%  '
# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08093'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.26%  '
# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.930.79'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.00%  '
# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.413'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.87%  '
# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '94.59'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -7.16%  '
# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.48'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.37%  '
# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.306.24'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.77%  '
# Row 18
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '255.44'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -7.49%  '
# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007992'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.43%  '
# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.809'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.34%  '
# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.185.22'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.53%  '
# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.13%  '
# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.003'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.32%  '
# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.854'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -5.46%  '
# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.666'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.85%  '
# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '159.12'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.93%  '
# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.361'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.99%  '
# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.09'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.40%  '
# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.1323'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -9.64%  '
# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.555'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.92%  '
# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.346'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.64%  '
# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.391'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.56%  '
# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.193'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.02%  '
# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05066'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.30%  '
# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.215'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.39%  '
# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7388'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.56%  '
# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.759'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.15%  '
# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01966'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.91%  '
# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.826'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.92%  '
# Row 40
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.582'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.34%  '
# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '78.99'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.76%  '
# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4444'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -7.24%  '
# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.991'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -9.35%  '
# Row 44
$ws.Range('E44').Value = '  +0.15%  '
# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8345'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.49%  '
# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.67'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.00%  '
# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.735'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.96%  '
# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.270'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -5.73%  '
# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '36.39'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.95%  '
# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05938'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.68%  '
# Row 51
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.467'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.88%  '
